$wb = $excel.ActiveWorkbook

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 26
$ws.Range("H26").Value = 5499.4
$ws.Range("I26").Value = 2749.5
$ws.Range("J26").Value = 7332.6665
$ws.Range("K26").Value = 2749.5
$ws.Range("L26").Value = 7332.6665
$ws.Range("M26").Value = -2419.5
$ws.Range("N26").Value = -7992.6665

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 82
$ws.Range("H82").Value = 10470.4
$ws.Range("I82").Value = 10470.4
$ws.Range("K82").Value = 10470.4
$ws.Range("M82").Value = -10087.4
# Row 85
$ws.Range("H85").Value = 10470.4
$ws.Range("I85").Value = 10470.4
$ws.Range("K85").Value = 10470.4
$ws.Range("M85").Value = -9144.4

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 3
$ws.Range("H3").Value = 4000
$ws.Range("I3").Value = 5000
$ws.Range("J3").Value = 3000
$ws.Range("K3").Value = 5000
$ws.Range("L3").Value = 3000
$ws.Range("M3").Value = -4887
$ws.Range("N3").Value = -3226

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 7
$ws.Range("H7").Value = 5500.3335
$ws.Range("I7").Value = 4250.5
$ws.Range("K7").Value = 4250.5
$ws.Range("M7").Value = -4138.5
# Row 8
$ws.Range("H8").Value = 5500.3335
$ws.Range("I8").Value = 4250.5
$ws.Range("K8").Value = 4250.5
$ws.Range("M8").Value = -4111.5
# Row 24
$ws.Range("H24").Value = 29999
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").Value = ""

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 3
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").Value = ""
# Row 5
$ws.Range("H5").Value = 15555
$ws.Range("I5").Value = 15555
$ws.Range("K5").Value = 15555
$ws.Range("M5").Value = -15442
# Row 14
$ws.Range("H14").Value = 1000
$ws.Range("I14").Value = 1000
$ws.Range("K14").Value = 1000
$ws.Range("M14").Value = -828
# Row 15
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").Value = ""
# Row 20
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").Value = ""
# Row 24
$ws.Range("H24").Value = 8000
$ws.Range("I24").Value = 8000
$ws.Range("K24").Value = 8000
$ws.Range("M24").Value = -7657
# Row 42
$ws.Range("H42").Value = 20333.334
$ws.Range("I42").Value = 14500
$ws.Range("J42").Value = 32000
$ws.Range("K42").Value = 14500
$ws.Range("L42").Value = 32000
$ws.Range("M42").Value = -13937
$ws.Range("N42").Value = -33126
# Row 43
$ws.Range("H43").Value = 9000
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").Value = ""
# Row 49
$ws.Range("H49").Value = 20333.334
$ws.Range("I49").Value = 14500
$ws.Range("J49").Value = 32000
$ws.Range("K49").Value = 14500
$ws.Range("L49").Value = 32000
$ws.Range("M49").Value = -14353
$ws.Range("N49").Value = -32294
# Row 50
$ws.Range("H50").Value = 25000
$ws.Range("I50").Value = 25000
$ws.Range("K50").Value = 25000
$ws.Range("M50").Value = -24363
# Row 56
$ws.Range("H56").Value = 10000
$ws.Range("I56").Value = 10000
$ws.Range("K56").Value = 10000
$ws.Range("M56").Value = -9309

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 3
$ws.Range("H3").Value = 1258888.2
$ws.Range("I3").Value = 2507777.5
$ws.Range("J3").Value = 9999
$ws.Range("K3").Value = 2507777.5
$ws.Range("L3").Value = 9999
$ws.Range("M3").Value = -2507663.5
$ws.Range("N3").Value = -10227
# Row 11
$ws.Range("H11").Value = 9999
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 9999
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 9999
$ws.Range("M11").Value = ""
$ws.Range("N11").Value = -10283
# Row 15
$ws.Range("H15").Value = 9999
$ws.Range("J15").Value = 9999
$ws.Range("L15").Value = 9999
$ws.Range("N15").Value = -10575
# Row 18
$ws.Range("H18").Value = 9999
$ws.Range("J18").Value = 9999
$ws.Range("L18").Value = 9999
$ws.Range("N18").Value = -10345
# Row 20
$ws.Range("H20").Value = 29999
$ws.Range("J20").Value = 29999
$ws.Range("L20").Value = 29999
$ws.Range("N20").Value = -30479
# Row 21
$ws.Range("H21").Value = 2513333.2
$ws.Range("J21").Value = 20000
$ws.Range("L21").Value = 20000
$ws.Range("N21").Value = -20470
# Row 22
$ws.Range("H22").Value = 29999
$ws.Range("J22").Value = 29999
$ws.Range("L22").Value = 29999
$ws.Range("N22").Value = -30585
# Row 24
$ws.Range("H24").Value = 1686666
$ws.Range("I24").Value = 5000000
$ws.Range("J24").Value = 29999
$ws.Range("K24").Value = 5000000
$ws.Range("L24").Value = 29999
$ws.Range("M24").Value = -4999770
$ws.Range("N24").Value = -30459
# Row 28
$ws.Range("H28").Value = 20000
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 20000
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 20000
$ws.Range("M28").Value = ""
$ws.Range("N28").Value = -20696
# Row 30
$ws.Range("H30").Value = 8266
$ws.Range("I30").Value = 4800
$ws.Range("J30").Value = 9999
$ws.Range("K30").Value = 4800
$ws.Range("L30").Value = 9999
$ws.Range("M30").Value = -4693
$ws.Range("N30").Value = -10213
# Row 31
$ws.Range("H31").Value = 14800
$ws.Range("I31").Value = 4400
$ws.Range("K31").Value = 4400
$ws.Range("M31").Value = -4052
# Row 33
$ws.Range("H33").Value = 5500
$ws.Range("I33").Value = 5500
$ws.Range("K33").Value = 5500
$ws.Range("M33").Value = -5250
# Row 35
$ws.Range("H35").Value = 2513333.2
$ws.Range("J35").Value = 20000
$ws.Range("L35").Value = 20000
$ws.Range("N35").Value = -20580
# Row 36
$ws.Range("H36").Value = 5500
$ws.Range("I36").Value = 5500
$ws.Range("K36").Value = 5500
$ws.Range("M36").Value = -5250
# Row 40
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").Value = ""
# Row 51
$ws.Range("H51").Value = 18000
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").Value = ""
# Row 52
$ws.Range("H52").Value = 10024500
$ws.Range("I52").Value = 20000000
$ws.Range("K52").Value = 20000000
$ws.Range("M52").Value = -19999774
# Row 58
$ws.Range("H58").Value = 8085
$ws.Range("I58").Value = 8085
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 8085
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -7777
$ws.Range("N58").Value = ""
# Row 135
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").Value = ""
# Row 138
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").Value = ""
